# Preços da Gasolina.xlsx — apply the "Update Projeto integrador" edit.
#
# Summary of the change:
#   1. The "gasolina_2010" table's PREÇO MÍNIMO / PREÇO MÁXIMO columns (F/G)
#      are renamed to PREÇO UNITARIO / PREÇO FINAL.
#   2. Column G becomes a calculated table column:
#         PREÇO UNITARIO * QUANTIDADE LTS
#      for every existing data row.
#   3. Column A (ID) is renumbered from the old external-query ids to a
#      simple sequential index (1, 2, 3, …) matching the row's position in
#      the table.
#   4. A new 60th data row (worksheet row 61) is appended with its own id,
#      region/state/product and price data, using the same formula.
#   5. The active selection moves to H61.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Meus Preços")
$lo = $ws.ListObjects.Item("gasolina_2010")

$lastRow = 60          # last pre-existing data row (worksheet row 2..60)
$newRow  = $lastRow + 1 # worksheet row 61 — the brand-new record

# ---------------------------------------------------------------------
# 1. Rename the table's price columns (editing the header cell drives the
#    ListObject's ListColumn.Name along with it).
# ---------------------------------------------------------------------
$ws.Range("F1").Value2 = "PREÇO UNITARIO"
$ws.Range("G1").Value2 = "PREÇO FINAL"

# ---------------------------------------------------------------------
# 2 & 3. Renumber the ID column and turn column G into the
#    PREÇO UNITARIO * QUANTIDADE LTS calculated column for every
#    pre-existing row.
# ---------------------------------------------------------------------
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 1).Value2 = $r - 1
    $ws.Range("G$r").Formula = "=gasolina_2010[[#This Row],[PREÇO UNITARIO]]*gasolina_2010[[#This Row],[QUANTIDADE LTS]]"
}

# ---------------------------------------------------------------------
# 4. Append the new 60th record (worksheet row 61).
# ---------------------------------------------------------------------
$ws.Cells.Item($newRow, 1).Value2 = $lastRow          # ID = 60
$ws.Cells.Item($newRow, 2).Value2 = "SUDESTE"          # REGIÃO
$ws.Cells.Item($newRow, 3).Value2 = "SAO PAULO"        # ESTADO
$ws.Cells.Item($newRow, 4).Value2 = "DIESEL COMUM"     # PRODUTO
$ws.Cells.Item($newRow, 5).Value2 = 3                  # QUANTIDADE LTS
$ws.Cells.Item($newRow, 6).Value2 = 2.99               # PREÇO UNITARIO
$ws.Range("G$newRow").Formula = "=gasolina_2010[[#This Row],[PREÇO UNITARIO]]*gasolina_2010[[#This Row],[QUANTIDADE LTS]]"

# ---------------------------------------------------------------------
# 5. Move the selection to H61, matching the saved view state.
# ---------------------------------------------------------------------
[void]$ws.Range("H61").Select()
